$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item("TextBox 3")
$tr = $sh.TextFrame.TextRange

# Original text: "Followed by a picture"
# Original runs: "Followed " | "by " | "a " | "picture"
# Target runs:   "Followed" | " " | "by" | " " | "a" | " " | "picture"
#
# Re-assign the text of each sub-range (by character position) to the same
# characters it already holds; PowerPoint's COM text-range model splits the
# underlying run wherever a TextRange.Characters() sub-range boundary does
# not coincide with an existing run boundary, and otherwise leaves the run
# untouched.

$tr.Characters(1, 8).Text = "Followed"
$tr.Characters(9, 1).Text = " "
$tr.Characters(10, 2).Text = "by"
$tr.Characters(12, 1).Text = " "
$tr.Characters(13, 1).Text = "a"
$tr.Characters(14, 1).Text = " "
